$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.477.66'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.375.54'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '506.73'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.55'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.32%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.387.25'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.14%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.80%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'Toncoin'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.85'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.78%  '
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Cardano'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.331'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.799.38'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.444.54'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.62'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.431.91'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.74%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '309.56'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.29'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.25%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.21'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.371'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.149'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.26'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.54'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0713'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.55%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.59%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.60%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.996'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.07'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.66'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.78%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.70'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.826'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.46'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.71%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '126.79'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.76'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.15%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.71%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '239.79'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0482'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.33%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.01'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.951'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.24%  '
